$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-preserving writes for the "Price" column (D) so that
# decimal-looking strings (e.g. "1.00", "0.100", "58.10") keep their
# exact textual representation instead of being coerced to numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '41.346.68'
$ws.Range("E2").Value = '  -5.93%  '
Set-TextValue $ws.Range("D3") '2.219.70'
$ws.Range("E3").Value = '  -6.22%  '
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue $ws.Range("D5") '243.95'
$ws.Range("E5").Value = '  +1.62%  '
Set-TextValue $ws.Range("D6") '0.620'
$ws.Range("E6").Value = '  -7.80%  '
Set-TextValue $ws.Range("D7") '69.66'
$ws.Range("E7").Value = '  -6.55%  '
$ws.Range("E8").Value = '  +0.16%  '
Set-TextValue $ws.Range("D9") '0.557'
$ws.Range("E9").Value = '  -7.63%  '
Set-TextValue $ws.Range("D10") '39.56'
$ws.Range("E10").Value = '  +5.81%  '
Set-TextValue $ws.Range("D11") '0.0955'
$ws.Range("E11").Value = '  -7.14%  '
Set-TextValue $ws.Range("D12") '58.10'
$ws.Range("E12").Value = '  -3.60%  '
$ws.Range("E13").Value = '  -3.84%  '
$ws.Range("E14").Value = '  -7.59%  '
Set-TextValue $ws.Range("D15") '2.551.08'
$ws.Range("E15").Value = '  -6.09%  '
Set-TextValue $ws.Range("D16") '14.78'
$ws.Range("E16").Value = '  -10.21%  '
Set-TextValue $ws.Range("D17") '0.843'
$ws.Range("E17").Value = '  -9.65%  '
Set-TextValue $ws.Range("D18") '2.218.90'
$ws.Range("E18").Value = '  -6.10%  '
Set-TextValue $ws.Range("D19") '41.296.44'
$ws.Range("E19").Value = '  -5.81%  '
Set-TextValue $ws.Range("D20") '0.0₃0949'
$ws.Range("E20").Value = '  -8.96%  '
Set-TextValue $ws.Range("D21") '72.25'
$ws.Range("E21").Value = '  -6.70%  '
$ws.Range("E22").Value = '  -8.27%  '
Set-TextValue $ws.Range("D23") '232.03'
$ws.Range("E23").Value = '  -8.79%  '
Set-TextValue $ws.Range("D24") '2.10'
$ws.Range("E24").Value = '  +11.92%  '
$ws.Range("E25").Value = '  +0.17%  '
Set-TextValue $ws.Range("D26") '3.62'
$ws.Range("E26").Value = '  -4.69%  '
Set-TextValue $ws.Range("D27") '2.42'
$ws.Range("E27").Value = '  -3.38%  '
$ws.Range("E28").Value = '  -7.71%  '
$ws.Range("E29").Value = '  -4.94%  '
Set-TextValue $ws.Range("D30") '172.25'
$ws.Range("E30").Value = '  -1.72%  '
Set-TextValue $ws.Range("D31") '20.49'
$ws.Range("E31").Value = '  -8.59%  '
$ws.Range("E32").Value = '  -8.35%  '
$ws.Range("E33").Value = '  -7.68%  '
Set-TextValue $ws.Range("D34") '0.0716'
$ws.Range("E34").Value = '  -6.17%  '
Set-TextValue $ws.Range("D35") '5.25'
$ws.Range("E35").Value = '  -4.17%  '
Set-TextValue $ws.Range("D36") '4.61'
$ws.Range("E36").Value = '  -10.14%  '
Set-TextValue $ws.Range("D37") '3.89'
$ws.Range("E37").Value = '  +2.09%  '
Set-TextValue $ws.Range("D38") '24.31'
$ws.Range("E38").Value = '  +17.08%  '
Set-TextValue $ws.Range("D39") '0.0277'
$ws.Range("E39").Value = '  -1.01%  '
$ws.Range("E40").Value = '  -5.09%  '
$ws.Range("E41").Value = '  -11.68%  '
Set-TextValue $ws.Range("D42") '65.95'
$ws.Range("E42").Value = '  +1.32%  '
Set-TextValue $ws.Range("D43") '4.99'
$ws.Range("E43").Value = '  -11.53%  '
Set-TextValue $ws.Range("D44") '0.203'
$ws.Range("E44").Value = '  +0.61%  '
Set-TextValue $ws.Range("D45") '8.82'
$ws.Range("E45").Value = '  -3.16%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D46") '0.100'
$ws.Range("E46").Value = '  -6.56%  '
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range("D47") '10.80'
$ws.Range("E47").Value = '  +11.52%  '
$ws.Range("B48").Value = 'SynthetixNetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextValue $ws.Range("D48") '4.63'
$ws.Range("E48").Value = '  +6.00%  '
$ws.Range("B49").Value = 'BinanceUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D49") '1.00'
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("E50").Value = '  -6.00%  '
$ws.Range("E51").Value = '  -5.67%  '
